$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) updates ---
$ws.Range("D2").Value = "'246.59"
$ws.Range("D3").Value = "'26.53"
$ws.Range("D5").Value = "'0.05613"
$ws.Range("D6").Value = "'6.479"
$ws.Range("D7").Value = "'0.8139"
$ws.Range("D8").Value = "'0.8453"
$ws.Range("D10").Value = "'0.02857"
$ws.Range("D11").Value = "'0.09390"
$ws.Range("D12").Value = "'0.001526"
$ws.Range("D19").Value = "'0.06966"
$ws.Range("D20").Value = "'0.03155"
$ws.Range("D22").Value = "'3.758"
$ws.Range("D23").Value = "'0.04667"
$ws.Range("D25").Value = "'0.001248"
$ws.Range("D26").Value = "'0.004589"
$ws.Range("D27").Value = "'0.00009606"
$ws.Range("D40").Value = "'0.03667"
$ws.Range("D44").Value = "'0.008908"
$ws.Range("D45").Value = "'0.00005295"
$ws.Range("D48").Value = "'0.002518"

# --- Row content changes (coin swapped into this row) ---
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.0006003"
$ws.Range("E9").Value = "8OneONE"

$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006167"
$ws.Range("E13").Value = "12TigerCashTCH"

$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.589"
$ws.Range("E14").Value = "13LEOLEO"

$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.011"
$ws.Range("E15").Value = "14GateTokenGT"

$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.118"
$ws.Range("E16").Value = "15BTSETokenBTSE"

$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3184"
$ws.Range("E17").Value = "16BitpandaEcosystemTokenBEST"

$ws.Range("B18").Value = "WazirX"
$ws.Range("C18").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D18").Value = "'0.1340"
$ws.Range("E18").Value = "17WazirXWRX"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006184"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"
